$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring over an (empty) cell from the row above so that the new row's
# Definition column (E) gets a real, present-but-blank cell, matching the
# existing "Definition" column pattern used throughout the sheet.
$ws.Range("E22").Copy($ws.Range("E23"))

# Columns A (Code) and C (Parent_Code) hold text-typed codes (e.g. "01",
# "02", "10" ...) even when the code looks numeric. Force text formatting
# before assigning so "99" / "9" are stored as text, consistent with the
# rest of the Code / Parent_Code columns.
$ws.Range("A23").NumberFormat = "@"
$ws.Range("C23").NumberFormat = "@"

$ws.Range("A23").Value = "99"
$ws.Range("B23").Value = "Not elsewhere classified"
$ws.Range("C23").Value = "9"
$ws.Range("D23").Value = "Not elsewhere classified"
